$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'313.14"
$ws.Range("D2").Style = $__style
$__style = $ws.Range("E2").Style
$ws.Range("E2").Value = "'2.53%"
$ws.Range("E2").Style = $__style
$__style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'37.79"
$ws.Range("D3").Style = $__style
$__style = $ws.Range("E3").Style
$ws.Range("E3").Value = "'1.72%"
$ws.Range("E3").Style = $__style
$__style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'5.148"
$ws.Range("D4").Style = $__style
$__style = $ws.Range("E4").Style
$ws.Range("E4").Value = "'0.85%"
$ws.Range("E4").Style = $__style
$__style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'0.07900"
$ws.Range("D5").Style = $__style
$__style = $ws.Range("E5").Style
$ws.Range("E5").Value = "'2.27%"
$ws.Range("E5").Style = $__style
$__style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'4.427"
$ws.Range("D6").Style = $__style
$__style = $ws.Range("E6").Style
$ws.Range("E6").Value = "'1.04%"
$ws.Range("E6").Style = $__style
$__style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'1.919"
$ws.Range("D7").Style = $__style
$__style = $ws.Range("E7").Style
$ws.Range("E7").Value = "'1.48%"
$ws.Range("E7").Style = $__style
$__style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'8.316"
$ws.Range("D8").Style = $__style
$__style = $ws.Range("E8").Style
$ws.Range("E8").Value = "'1.25%"
$ws.Range("E8").Style = $__style
$__style = $ws.Range("E9").Style
$ws.Range("E9").Value = "'-9.78%"
$ws.Range("E9").Style = $__style
$__style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.9216"
$ws.Range("D10").Style = $__style
$__style = $ws.Range("E10").Style
$ws.Range("E10").Value = "'0.31%"
$ws.Range("E10").Style = $__style
$__style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.1239"
$ws.Range("D11").Style = $__style
$__style = $ws.Range("E11").Style
$ws.Range("E11").Value = "'0.06%"
$ws.Range("E11").Style = $__style
$__style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.1923"
$ws.Range("D12").Style = $__style
$__style = $ws.Range("E12").Style
$ws.Range("E12").Value = "'2.38%"
$ws.Range("E12").Style = $__style
$__style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.09165"
$ws.Range("D13").Style = $__style
$__style = $ws.Range("E13").Style
$ws.Range("E13").Value = "'4.62%"
$ws.Range("E13").Style = $__style
$__style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.03345"
$ws.Range("D14").Style = $__style
$__style = $ws.Range("E14").Style
$ws.Range("E14").Value = "'-2.50%"
$ws.Range("E14").Style = $__style
$__style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.09614"
$ws.Range("D15").Style = $__style
$__style = $ws.Range("E15").Style
$ws.Range("E15").Value = "'-1.03%"
$ws.Range("E15").Style = $__style
$__style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'0.001380"
$ws.Range("D16").Style = $__style
$__style = $ws.Range("E16").Style
$ws.Range("E16").Value = "'0.62%"
$ws.Range("E16").Style = $__style
$__style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.005729"
$ws.Range("D17").Style = $__style
$__style = $ws.Range("E17").Style
$ws.Range("E17").Value = "'-7.19%"
$ws.Range("E17").Style = $__style
$__style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'3.499"
$ws.Range("D18").Style = $__style
$__style = $ws.Range("E18").Style
$ws.Range("E18").Value = "'-1.66%"
$ws.Range("E18").Style = $__style
$__style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.3443"
$ws.Range("D19").Style = $__style
$__style = $ws.Range("E19").Style
$ws.Range("E19").Value = "'2.09%"
$ws.Range("E19").Style = $__style
$__style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'5.281"
$ws.Range("D20").Style = $__style
$__style = $ws.Range("E20").Style
$ws.Range("E20").Value = "'4.93%"
$ws.Range("E20").Style = $__style
$__style = $ws.Range("E21").Style
$ws.Range("E21").Value = "'-0.79%"
$ws.Range("E21").Style = $__style
$__style = $ws.Range("E22").Style
$ws.Range("E22").Value = "'3.52%"
$ws.Range("E22").Style = $__style
$__style = $ws.Range("E23").Style
$ws.Range("E23").Value = "'-0.52%"
$ws.Range("E23").Style = $__style
$__style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'0.04365"
$ws.Range("D24").Style = $__style
$__style = $ws.Range("E24").Style
$ws.Range("E24").Value = "'0.72%"
$ws.Range("E24").Style = $__style
$__style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.001250"
$ws.Range("D25").Style = $__style
$__style = $ws.Range("E25").Style
$ws.Range("E25").Value = "'2.49%"
$ws.Range("E25").Style = $__style
$__style = $ws.Range("E26").Style
$ws.Range("E26").Value = "'-3.28%"
$ws.Range("E26").Style = $__style
$__style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'0.0001221"
$ws.Range("D27").Style = $__style
$__style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.02243"
$ws.Range("D39").Style = $__style
$__style = $ws.Range("E39").Style
$ws.Range("E39").Value = "'1.58%"
$ws.Range("E39").Style = $__style
$__style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.05105"
$ws.Range("D40").Style = $__style
$__style = $ws.Range("E40").Style
$ws.Range("E40").Value = "'3.97%"
$ws.Range("E40").Style = $__style
$__style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.007477"
$ws.Range("D41").Style = $__style
$__style = $ws.Range("E41").Style
$ws.Range("E41").Value = "'-2.02%"
$ws.Range("E41").Style = $__style
$__style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.1364"
$ws.Range("D42").Style = $__style
$__style = $ws.Range("E42").Style
$ws.Range("E42").Value = "'2.37%"
$ws.Range("E42").Style = $__style
$__style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.008805"
$ws.Range("D43").Style = $__style
$__style = $ws.Range("E43").Style
$ws.Range("E43").Value = "'-11.20%"
$ws.Range("E43").Style = $__style
$__style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.001922"
$ws.Range("D44").Style = $__style
$__style = $ws.Range("E44").Style
$ws.Range("E44").Value = "'-4.07%"
$ws.Range("E44").Style = $__style
$__style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.008628"
$ws.Range("D45").Style = $__style
$__style = $ws.Range("E45").Style
$ws.Range("E45").Value = "'-2.03%"
$ws.Range("E45").Style = $__style
$__style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.00006739"
$ws.Range("D46").Style = $__style
$__style = $ws.Range("E46").Style
$ws.Range("E46").Value = "'-3.22%"
$ws.Range("E46").Style = $__style
$__style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = $__style
$__style = $ws.Range("E47").Style
$ws.Range("E47").Value = "'-0.44%"
$ws.Range("E47").Style = $__style
$__style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.003352"
$ws.Range("D48").Style = $__style
$__style = $ws.Range("E48").Style
$ws.Range("E48").Value = "'11.28%"
$ws.Range("E48").Style = $__style
$__style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.001201"
$ws.Range("D49").Style = $__style
$__style = $ws.Range("E49").Style
$ws.Range("E49").Value = "'-8.12%"
$ws.Range("E49").Style = $__style
$__style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = $__style
$__style = $ws.Range("E50").Style
$ws.Range("E50").Value = "'-0.44%"
$ws.Range("E50").Style = $__style
$__style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = $__style
$__style = $ws.Range("E51").Style
$ws.Range("E51").Value = "'-0.44%"
$ws.Range("E51").Style = $__style
